$d = $word.ActiveDocument

# Bold the entire first paragraph (including paragraph mark)
$p1 = $d.Paragraphs(1)
$p1.Range.Bold = 1
$p1.Range.BoldBi = 1
